$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT, even when it looks numeric
# (the sheet stores "Price" figures as inline strings, not numbers).
# Briefly forcing a text number format prevents Excel from auto-converting
# the digits to a real number, then ClearFormats() drops that temporary
# format again so the cell's style is left untouched.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Updated "Price" (column D) readings
$prices = [ordered]@{
    "D2"  = "266.72"
    "D3"  = "21.54"
    "D4"  = "6.173"
    "D5"  = "0.06164"
    "D6"  = "3.570"
    "D7"  = "6.531"
    "D8"  = "1.361"
    "D9"  = "0.8226"
    "D10" = "0.01343"
    "D11" = "0.1561"
    "D12" = "0.08187"
    "D13" = "0.03338"
    "D14" = "0.03174"
    "D15" = "0.09255"
    "D16" = "3.760"
    "D17" = "0.001624"
    "D18" = "0.04675"
    "D19" = "0.006374"
    "D20" = "0.006203"
    "D21" = "0.001068"
    "D22" = "0.0001499"
    "D25" = "0.3301"
    "D28" = "0.0001617"
    "D40" = "0.04638"
    "D41" = "0.006971"
    "D42" = "0.1128"
    "D43" = "0.003658"
    "D45" = "0.00005921"
    "D46" = "0.0009892"
    "D48" = "0.7815"
    "D49" = "0.002440"
    "D50" = "0.00001899"
    "D51" = "0.01239"
}

foreach ($addr in $prices.Keys) {
    Set-TextValue $addr $prices[$addr]
}

# Rows 42 and 43 also swapped places in the symbol list: row 42 is now
# BKEXToken and row 43 is now CEJI (their Price values were updated above).
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
